# Refresh the crypto price/volume snapshot (columns D and E, rows 2-51).
#
# The Price column (D) holds plain text in the workbook (t="inlineStr"),
# e.g. "1.008" or "26.228.27" -- not real numbers. Assigning such a
# numeric-looking string straight to Range.Value makes Excel parse it as a
# number (losing exact digits, e.g. trailing zeros) and also flips the
# cell's style (quote-prefix / text number format) versus the untouched
# style-less original cell. To keep both the literal text AND the original
# (default) cell style, such values are staged in a scratch cell that is
# explicitly formatted as text, then brought over with a values-only paste
# (xlPasteValues), which copies the text content but not the formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$scratch = $ws.Range("Z1")

$ws.Range("D2").Value = '26.228.27'
$ws.Range("E2").Value = '  +0.49%  '

$ws.Range("D3").Value = '1.658.65'
$ws.Range("E3").Value = '  +0.09%  '

$scratch.NumberFormat = "@"
$scratch.Value = '1.008'
$scratch.Copy()
$ws.Range("D4").PasteSpecial(-4163)
$ws.Range("E4").Value = '  +0.63%  '

$scratch.NumberFormat = "@"
$scratch.Value = '218.07'
$scratch.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = '  -0.20%  '

$scratch.NumberFormat = "@"
$scratch.Value = '0.5305'
$scratch.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = '  +0.09%  '

$scratch.NumberFormat = "@"
$scratch.Value = '1.009'
$scratch.Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("E7").Value = '  +0.62%  '

$scratch.NumberFormat = "@"
$scratch.Value = '0.2632'
$scratch.Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value = '  +0.61%  '

$scratch.NumberFormat = "@"
$scratch.Value = '0.06351'
$scratch.Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = '  +0.20%  '

$scratch.NumberFormat = "@"
$scratch.Value = '20.44'
$scratch.Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = '  +0.01%  '

$scratch.NumberFormat = "@"
$scratch.Value = '0.07840'
$scratch.Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = '  +0.91%  '

$scratch.NumberFormat = "@"
$scratch.Value = '4.541'
$scratch.Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = '  +0.91%  '

$ws.Range("D13").Value = '1.663.29'
$ws.Range("E13").Value = '  +0.37%  '

$ws.Range("D14").Value = '1.889.94'
$ws.Range("E14").Value = '  +0.32%  '

$scratch.NumberFormat = "@"
$scratch.Value = '0.5525'
$scratch.Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value = '  +0.86%  '

$ws.Range("D16").Value = '0.0₅8156'
$ws.Range("E16").Value = '  -0.14%  '

$scratch.NumberFormat = "@"
$scratch.Value = '65.61'
$scratch.Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("E17").Value = '  +0.54%  '

$ws.Range("D18").Value = '26.234.40'
$ws.Range("E18").Value = '  +0.40%  '

$scratch.NumberFormat = "@"
$scratch.Value = '1.008'
$scratch.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = '  +0.62%  '

$scratch.NumberFormat = "@"
$scratch.Value = '4.643'
$scratch.Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = '  +2.08%  '

$scratch.NumberFormat = "@"
$scratch.Value = '192.69'
$scratch.Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = '  -0.46%  '

$scratch.NumberFormat = "@"
$scratch.Value = '10.17'
$scratch.Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = '  +0.85%  '

$scratch.NumberFormat = "@"
$scratch.Value = '6.037'
$scratch.Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = '  +0.22%  '

$scratch.NumberFormat = "@"
$scratch.Value = '1.010'
$scratch.Copy()
$ws.Range("D24").PasteSpecial(-4163)

$scratch.NumberFormat = "@"
$scratch.Value = '144.38'
$scratch.Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = '  +2.84%  '

$scratch.NumberFormat = "@"
$scratch.Value = '0.1220'
$scratch.Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = '  -2.06%  '

$scratch.NumberFormat = "@"
$scratch.Value = '7.214'
$scratch.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = '  -0.90%  '

$scratch.NumberFormat = "@"
$scratch.Value = '16.16'
$scratch.Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = '  -0.04%  '

$ws.Range("E29").Value = '  +3.24%  '

$scratch.NumberFormat = "@"
$scratch.Value = '0.05850'
$scratch.Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = '  -1.67%  '

$scratch.NumberFormat = "@"
$scratch.Value = '1.277'
$scratch.Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = '  -0.07%  '

$scratch.NumberFormat = "@"
$scratch.Value = '3.572'
$scratch.Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = '  +1.74%  '

$scratch.NumberFormat = "@"
$scratch.Value = '3.288'
$scratch.Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = '  +1.44%  '

$scratch.NumberFormat = "@"
$scratch.Value = '1.612'
$scratch.Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = '  +3.53%  '

$scratch.NumberFormat = "@"
$scratch.Value = '0.9581'
$scratch.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = '  +0.94%  '

$scratch.NumberFormat = "@"
$scratch.Value = '2.817'
$scratch.Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = '  +1.66%  '

$scratch.NumberFormat = "@"
$scratch.Value = '2.421'
$scratch.Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = '  +0.40%  '

$scratch.NumberFormat = "@"
$scratch.Value = '0.5801'
$scratch.Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = '  +2.79%  '

$scratch.NumberFormat = "@"
$scratch.Value = '0.01606'
$scratch.Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = '  -0.45%  '

$scratch.NumberFormat = "@"
$scratch.Value = '5.887'
$scratch.Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = '  +0.67%  '

$scratch.NumberFormat = "@"
$scratch.Value = '0.8511'
$scratch.Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = '  +0.31%  '

$scratch.NumberFormat = "@"
$scratch.Value = '1.008'
$scratch.Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = '  +0.59%  '

$ws.Range("D43").Value = '1.042.93'
$ws.Range("E43").Value = '  +2.96%  '

$scratch.NumberFormat = "@"
$scratch.Value = '103.87'
$scratch.Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = '  +2.39%  '

$ws.Range("D45").Value = '1.803.56'
$ws.Range("E45").Value = '  +0.21%  '

$scratch.NumberFormat = "@"
$scratch.Value = '57.08'
$scratch.Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = '  -0.07%  '

$ws.Range("D47").Value = '0.0₈108'
$ws.Range("E47").Value = '  +3.36%  '

$scratch.NumberFormat = "@"
$scratch.Value = '1.012'
$scratch.Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = '  +0.92%  '

$scratch.NumberFormat = "@"
$scratch.Value = '0.4369'
$scratch.Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = '  +1.95%  '

$scratch.NumberFormat = "@"
$scratch.Value = '7.964'
$scratch.Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = '  +2.92%  '

$scratch.NumberFormat = "@"
$scratch.Value = '0.05161'
$scratch.Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = '  +0.10%  '

# Clean up the scratch cell so it leaves no trace in the sheet.
$scratch.Clear()
$excel.CutCopyMode = $false
